# Commit: "add the NA's under duplicate_image_filename"
# Column E ("duplicate_image_filename", see header E1) was only populated
# for the header row and the summary table at the bottom (E28). Every data
# row in the main stimuli table (rows 2-21) was missing a value for that
# column - fill them in with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
